$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22 gains a brand-new results record (model #14), mirroring the shape
# of the preceding rows (row 21, model #13) in the "Predicting N samples
# using past 10" table.

# Seed the formatting for the new row by cloning row 21's formats first,
# so the new cells pick up the same style ("s=4") used throughout the table.
$ws.Range("A21:K21").Copy() | Out-Null
$ws.Range("A22:K22").PasteSpecial(-4122) | Out-Null

# Numeric columns: #Model, Layers, Nodes, Dropout, Future samples (N)
$ws.Range("A22").Value = 14.0
$ws.Range("B22").Value = 1.0
$ws.Range("C22").Value = 50.0
$ws.Range("D22").Value = 0.0
$ws.Range("E22").Value = 5.0

# Text columns: leading apostrophe forces these to be stored as text
# (shared strings) instead of being auto-parsed as numbers/percentages.
$ws.Range("F22").Value = "'2.981210708618164"
$ws.Range("G22").Value = "'1m 19s / 1m 6s"
$ws.Range("H22").Value = "'12.592828902855075 %"
$ws.Range("I22").Value = "'5.987421332258472 %"
$ws.Range("J22").Value = "'17.05721783694347 %"
$ws.Range("K22").Value = "'6/5 epochs"

# Re-apply row 21's formatting so the quote-prefix/text number format that
# Value assignment above implicitly introduced doesn't stick around - the
# new row should keep the same plain style as the rest of the table.
$ws.Range("A21:K21").Copy() | Out-Null
$ws.Range("A22:K22").PasteSpecial(-4122) | Out-Null
